$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new Kaspa buy recorded on 2025-06-22 as row 35.
$row = 35

# Column A holds the date as literal text (e.g. "06/22/2025"), not an
# Excel date value, so force a text format before assigning it -- this
# matches how the other recently-appended rows (e.g. row 34) store it.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "06/22/2025"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 155.4730000000054
$ws.Cells.Item($row, 3).Value = 0.06431984974882875
$ws.Cells.Item($row, 4).Value = 10
